$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally has 11 columns:
#   A Time | B Teams | C Goalies | D Win | E Best ML | F Best Spread |
#   G Goals | H Total Goals | I Best O/U | J Bet Value | K More Details
#
# The edit drops the Goalies, Best ML, Best Spread, Bet Value and
# More Details columns entirely (their data, header strings, hyperlinks
# and related styles all go away), leaving:
#   A Time | B Teams | C Win | D Goals | E Total Goals | F Best O/U
#
# Delete from right-most column to left-most so earlier deletions don't
# shift the column letters of the ones still to be removed.
$ws.Columns("K").Delete()
$ws.Columns("J").Delete()
$ws.Columns("F").Delete()
$ws.Columns("E").Delete()
$ws.Columns("C").Delete()
